# Applies the "cryptos list" price/volume refresh described in the commit
# message (GitHub Actions scheduled update). Column D (Price) and column E
# (Volume(1h)) are plain text cells in the workbook (values like "2.300.29"
# or "0.0₃0908" are not valid numbers, and the "  +1.29%  " volume strings
# carry significant leading/trailing spaces) so every D/E write below is
# apostrophe-prefixed to force a literal text entry -- this stops Excel from
# re-interpreting number-looking text (e.g. "2.77", "0.100") as a numeric or
# date value. The Style reset that follows each write clears the transient
# "quote prefix" cell format Excel applies for that forced-text entry, so the
# cells keep the same (default/unstyled) look as their untouched neighbours.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($address, $text) {
    $ws.Range($address).Value = "'" + $text
    $ws.Range($address).Style = 'Normal'
}

Set-TextCell 'D2' '42.890.56'
Set-TextCell 'E2' '  -0.59%  '
Set-TextCell 'D3' '2.300.29'
Set-TextCell 'E3' '  -0.93%  '
Set-TextCell 'E4' '  +0.06%  '
Set-TextCell 'D5' '299.82'
Set-TextCell 'E5' '  -1.22%  '
Set-TextCell 'D6' '97.64'
Set-TextCell 'E6' '  -2.30%  '
Set-TextCell 'E7' '  +1.29%  '
Set-TextCell 'E8' '  +0.06%  '
Set-TextCell 'D9' '0.505'
Set-TextCell 'E9' '  -2.67%  '
Set-TextCell 'D10' '35.93'
Set-TextCell 'E10' '  -0.93%  '
Set-TextCell 'D11' '0.0787'
Set-TextCell 'E11' '  -0.53%  '
Set-TextCell 'E12' '  +0.75%  '
Set-TextCell 'D13' '17.73'
Set-TextCell 'E13' '  -0.34%  '
Set-TextCell 'D14' '6.76'
Set-TextCell 'E14' '  -2.18%  '
Set-TextCell 'D15' '2.659.04'
Set-TextCell 'E15' '  -0.80%  '
Set-TextCell 'D16' '2.285.89'
Set-TextCell 'E16' '  -1.12%  '
Set-TextCell 'D17' '0.778'
Set-TextCell 'E17' '  -2.31%  '
Set-TextCell 'D18' '42.880.46'
Set-TextCell 'E18' '  -0.39%  '
Set-TextCell 'D19' '12.57'
Set-TextCell 'E19' '  -4.06%  '
Set-TextCell 'D20' '0.0₃0908'
Set-TextCell 'E20' '  -0.37%  '
Set-TextCell 'D21' '6.07'
Set-TextCell 'E21' '  -2.07%  '
Set-TextCell 'D22' '67.97'
Set-TextCell 'E22' '  -0.47%  '
Set-TextCell 'D23' '242.77'
Set-TextCell 'E23' '  +1.09%  '
Set-TextCell 'D24' '2.13'
Set-TextCell 'E24' '  -1.67%  '
Set-TextCell 'E25' '  +0.03%  '
Set-TextCell 'D26' '2.42'
Set-TextCell 'E26' '  -1.86%  '
Set-TextCell 'E27' '  -0.29%  '
Set-TextCell 'D28' '25.19'
Set-TextCell 'E28' '  -1.55%  '
Set-TextCell 'D29' '165.95'
Set-TextCell 'E29' '  -1.88%  '
Set-TextCell 'D30' '2.03'
Set-TextCell 'E30' '  -1.02%  '
Set-TextCell 'D31' '9.03'
Set-TextCell 'E31' '  -1.85%  '
Set-TextCell 'D32' '32.75'
Set-TextCell 'E32' '  -4.18%  '
Set-TextCell 'E33' '  +0.11%  '
Set-TextCell 'D34' '4.81'
Set-TextCell 'E34' '  -3.27%  '
Set-TextCell 'D35' '4.99'
Set-TextCell 'E35' '  -3.64%  '
Set-TextCell 'D36' '17.19'
Set-TextCell 'E36' '  -3.96%  '
Set-TextCell 'E37' '  -0.58%  '
Set-TextCell 'D38' '0.0687'
Set-TextCell 'E38' '  -1.64%  '
Set-TextCell 'D39' '0.100'
Set-TextCell 'E39' '  -2.39%  '
Set-TextCell 'E40' '  -4.01%  '
Set-TextCell 'D41' '2.74'
Set-TextCell 'E41' '  -1.56%  '
Set-TextCell 'E42' '  +0.03%  '
Set-TextCell 'D43' '2.008.40'
Set-TextCell 'E43' '  +0.74%  '
Set-TextCell 'D44' '0.0284'
Set-TextCell 'E44' '  -1.74%  '
Set-TextCell 'D45' '10.13'
Set-TextCell 'E45' '  -0.57%  '
Set-TextCell 'D46' '2.13'
Set-TextCell 'E46' '  -5.27%  '
Set-TextCell 'D47' '17.28'
Set-TextCell 'E47' '  -2.02%  '
Set-TextCell 'D48' '2.77'
Set-TextCell 'E48' '  -3.20%  '
Set-TextCell 'D49' '2.524.59'
Set-TextCell 'E49' '  -0.85%  '
Set-TextCell 'D50' '53.08'
Set-TextCell 'E50' '  -3.58%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D51' '2.77'
Set-TextCell 'E51' '  -2.84%  '
